# Generate Report for Handback
# Update status text from "Ready for handoff" to "Handback transform failed"
# for the a72fb9c5 file, on the Overview sheet and each locale sheet, and
# record the handback/handoff file name mismatch error detail for each locale.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# Overview sheet: row for a72fb9c5-... has both the zh-cn and de-de status
# in columns B and C.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# zh-cn sheet: Status column (C) for the a72fb9c5 row, plus new Error Detail
# (column L) describing the handback/handoff file name mismatch.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("L3").Value = "Handback file name: au4stwhx.at4 is different with handoff file name: a72fb9c5-0acd-442d-b70d-1daaf0192f28.a81056532d2314408d4c4695976e23b6acd44b2c.zh-cn."

# de-de sheet: Status column (C) for the a72fb9c5 row, plus new Error Detail
# (column L) describing the handback/handoff file name mismatch.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("L3").Value = "Handback file name: au4stwhx.at4 is different with handoff file name: a72fb9c5-0acd-442d-b70d-1daaf0192f28.a81056532d2314408d4c4695976e23b6acd44b2c.de-de."
